# Update Leviathan_Profits market-price-derived cells (currentAveragePrice*, Leve*Price*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 21199.8
$ws.Range("I54").Value = 16666.666
$ws.Range("J54").Value = 27999.5
$ws.Range("K54").Value = 16666.666
$ws.Range("L54").Value = 27999.5
$ws.Range("M54").Value = -16180.666
$ws.Range("N54").Value = -28971.5
$ws.Range("H137").Value = 2582.375
$ws.Range("I137").Value = 2281.6086
$ws.Range("K137").Value = 6844.825800000001
$ws.Range("M137").Value = -4294.825800000001
$ws.Range("H138").Value = 1873.7693
$ws.Range("I138").Value = 1207.9706
$ws.Range("J138").Value = 2388.25
$ws.Range("K138").Value = 3623.9118
$ws.Range("L138").Value = 7164.75
$ws.Range("M138").Value = 1516.0882
$ws.Range("N138").Value = -17444.75
$ws.Range("H141").Value = 36933.258
$ws.Range("I141").Value = 36933.258
$ws.Range("K141").Value = 110799.774
$ws.Range("M141").Value = -105619.774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19324.812
$ws.Range("I32").Value = 3505.5386
$ws.Range("K32").Value = 3505.5386
$ws.Range("M32").Value = -3218.5386
$ws.Range("H45").Value = 724000.5600000001
$ws.Range("I45").Value = 1264388.1
$ws.Range("K45").Value = 1264388.1
$ws.Range("M45").Value = -1264011.1
$ws.Range("H74").Value = 1626.7
$ws.Range("I74").Value = 1346.2667
$ws.Range("K74").Value = 1346.2667
$ws.Range("M74").Value = -472.2666999999999
$ws.Range("H77").Value = 1626.7
$ws.Range("I77").Value = 1346.2667
$ws.Range("K77").Value = 6731.3335
$ws.Range("M77").Value = -2363.3335
$ws.Range("H132").Value = 3204.6667
$ws.Range("I132").Value = 3204.6667
$ws.Range("K132").Value = 9614.000100000001
$ws.Range("M132").Value = -7084.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1438.6825
$ws.Range("I134").Value = 1237.9323
$ws.Range("K134").Value = 3713.7969
$ws.Range("M134").Value = -1178.7969

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1186.4828
$ws.Range("I16").Value = 1193.6786
$ws.Range("K16").Value = 1193.6786
$ws.Range("M16").Value = -906.6786
$ws.Range("H31").Value = 20274.781
$ws.Range("I31").Value = 30217
$ws.Range("J31").Value = 2875.9
$ws.Range("K31").Value = 30217
$ws.Range("L31").Value = 2875.9
$ws.Range("M31").Value = -29922
$ws.Range("N31").Value = -3465.9
$ws.Range("H34").Value = 20274.781
$ws.Range("I34").Value = 30217
$ws.Range("J34").Value = 2875.9
$ws.Range("K34").Value = 30217
$ws.Range("L34").Value = 2875.9
$ws.Range("M34").Value = -30015
$ws.Range("N34").Value = -3279.9
$ws.Range("H58").Value = 1407
$ws.Range("I58").Value = 1071.4445
$ws.Range("J58").Value = 2413.6667
$ws.Range("K58").Value = 1071.4445
$ws.Range("L58").Value = 2413.6667
$ws.Range("M58").Value = -868.4445000000001
$ws.Range("N58").Value = -2819.6667
$ws.Range("H76").Value = 5250
$ws.Range("I76").Value = 5250
$ws.Range("K76").Value = 5250
$ws.Range("M76").Value = -4935
$ws.Range("H79").Value = 5250
$ws.Range("I79").Value = 5250
$ws.Range("K79").Value = 5250
$ws.Range("M79").Value = -4158
$ws.Range("H105").Value = 1536.9445
$ws.Range("I105").Value = 1616.5385
$ws.Range("K105").Value = 1616.5385
$ws.Range("M105").Value = 130.4614999999999
$ws.Range("H113").Value = 1186.4828
$ws.Range("I113").Value = 1193.6786
$ws.Range("K113").Value = 1193.6786
$ws.Range("M113").Value = 976.3214
$ws.Range("H132").Value = 2350.8542
$ws.Range("I132").Value = 2273.875
$ws.Range("K132").Value = 6821.625
$ws.Range("M132").Value = -4291.625
$ws.Range("H136").Value = 1407
$ws.Range("I136").Value = 1071.4445
$ws.Range("J136").Value = 2413.6667
$ws.Range("K136").Value = 3214.3335
$ws.Range("L136").Value = 7241.000100000001
$ws.Range("M136").Value = -664.3335000000002
$ws.Range("N136").Value = -12341.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1112.5
$ws.Range("I5").Value = 1167
$ws.Range("J5").Value = 1058
$ws.Range("K5").Value = 3501
$ws.Range("L5").Value = 3174
$ws.Range("M5").Value = -3389
$ws.Range("N5").Value = -3398
$ws.Range("H29").Value = 111154.336
$ws.Range("I29").Value = 166717.67
$ws.Range("J29").Value = 27.666666
$ws.Range("K29").Value = 500153.01
$ws.Range("L29").Value = 82.99999800000001
$ws.Range("M29").Value = -499876.01
$ws.Range("N29").Value = -636.999998
$ws.Range("H46").Value = 200662.8
$ws.Range("I46").Value = 200662.8
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 601988.3999999999
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -601897.3999999999
$ws.Range("N46").ClearContents()
$ws.Range("H128").Value = 499988
$ws.Range("I128").Value = 499988
$ws.Range("K128").Value = 1499964
$ws.Range("M128").Value = -1494984
$ws.Range("H131").Value = 45380.08
$ws.Range("I131").Value = 65209.117
$ws.Range("J131").Value = 3243.375
$ws.Range("K131").Value = 195627.351
$ws.Range("L131").Value = 9730.125
$ws.Range("M131").Value = -190587.351
$ws.Range("N131").Value = -19810.125
$ws.Range("H135").Value = 1112.5
$ws.Range("I135").Value = 1167
$ws.Range("J135").Value = 1058
$ws.Range("K135").Value = 10503
$ws.Range("L135").Value = 9522
$ws.Range("M135").Value = -7968
$ws.Range("N135").Value = -14592
$ws.Range("H140").Value = 3690.7058
$ws.Range("I140").Value = 3284.1667
$ws.Range("K140").Value = 9852.500100000001
$ws.Range("M140").Value = -4672.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 15880.75
$ws.Range("J98").Value = 15880.75
$ws.Range("L98").Value = 15880.75
$ws.Range("N98").Value = -21870.75
$ws.Range("H122").Value = 3601.8572
$ws.Range("I122").Value = 3601.8572
$ws.Range("K122").Value = 10805.5716
$ws.Range("M122").Value = -8355.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1683.3334
$ws.Range("I22").Value = 1533.3334
$ws.Range("J22").Value = 1833.3334
$ws.Range("K22").Value = 1533.3334
$ws.Range("L22").Value = 1833.3334
$ws.Range("M22").Value = -1238.3334
$ws.Range("N22").Value = -2423.3334
$ws.Range("H27").Value = 1683.3334
$ws.Range("I27").Value = 1533.3334
$ws.Range("J27").Value = 1833.3334
$ws.Range("K27").Value = 1533.3334
$ws.Range("L27").Value = 1833.3334
$ws.Range("M27").Value = -1426.3334
$ws.Range("N27").Value = -2047.3334
$ws.Range("H46").Value = 29738.334
$ws.Range("I46").Value = 71720.836
$ws.Range("K46").Value = 71720.836
$ws.Range("M46").Value = -71532.836
$ws.Range("H68").Value = 1817.5
$ws.Range("I68").Value = 1323.125
$ws.Range("K68").Value = 1323.125
$ws.Range("M68").Value = -574.125
$ws.Range("H71").Value = 1817.5
$ws.Range("I71").Value = 1323.125
$ws.Range("K71").Value = 6615.625
$ws.Range("M71").Value = -2871.625
$ws.Range("H132").Value = 2785.9492
$ws.Range("I132").Value = 2365.532
$ws.Range("J132").Value = 4432.5835
$ws.Range("K132").Value = 7096.596
$ws.Range("L132").Value = 13297.7505
$ws.Range("M132").Value = -4566.596
$ws.Range("N132").Value = -18357.7505
$ws.Range("H136").Value = 2148.0942
$ws.Range("I136").Value = 1791.5106
$ws.Range("K136").Value = 5374.531800000001
$ws.Range("M136").Value = -2824.531800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1118244
$ws.Range("I14").Value = 1675366.8
$ws.Range("J14").Value = 3998.5
$ws.Range("K14").Value = 1675366.8
$ws.Range("L14").Value = 3998.5
$ws.Range("M14").Value = -1675198.8
$ws.Range("N14").Value = -4334.5
$ws.Range("H107").Value = 100004350
$ws.Range("I107").Value = 6388.5
$ws.Range("K107").Value = 19165.5
$ws.Range("M107").Value = -17245.5
$ws.Range("H113").Value = 462.8
$ws.Range("I113").Value = 347.55554
$ws.Range("K113").Value = 1042.66662
$ws.Range("M113").Value = 1127.33338
$ws.Range("H122").Value = 2259.1365
$ws.Range("I122").Value = 2215.7896
$ws.Range("K122").Value = 6647.3688
$ws.Range("M122").Value = -4197.3688
$ws.Range("H132").Value = 2077979.1
$ws.Range("I132").Value = 1201320.1
$ws.Range("K132").Value = 3603960.3
$ws.Range("M132").Value = -3601430.3
